$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.114.13"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.651.52"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "218.10"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "0.5233"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").Value = "0.2614"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "0.06276"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").Value = "20.47"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "0.07821"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "4.470"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "1.659.31"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "1.878.94"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "0.5520"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "0.0₅7988"
$ws.Range("E16").Value = "  -2.75%  "
$ws.Range("D17").Value = "64.84"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "26.095.73"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "4.619"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").Value = "194.65"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "10.05"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "5.935"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "146.64"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "0.1203"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").Value = "7.148"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "15.90"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "1.483"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "0.05686"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("D31").Value = "1.268"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "3.466"
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("D33").Value = "3.323"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").Value = "1.580"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("D35").Value = "2.792"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").Value = "2.414"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "0.9464"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "0.5645"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "0.01591"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").Value = "5.917"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").Value = "1.060.92"
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").Value = "1.006"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").Value = "0.8418"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").Value = "103.15"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "1.793.32"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "57.36"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("D48").Value = "0.05403"
$ws.Range("E48").Value = "  +4.70%  "
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "0.4397"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "7.967"
$ws.Range("E51").Value = "  -0.87%  "

